# Auto-generated Excel COM-interop script to update column F (visitor/attendance counts)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 358
$ws.Cells.Item(3, 6).Value = 3604
$ws.Cells.Item(5, 6).Value = 8388
$ws.Cells.Item(7, 6).Value = 139
$ws.Cells.Item(8, 6).Value = 2278
$ws.Cells.Item(13, 6).Value = 7610
$ws.Cells.Item(14, 6).Value = 7824
$ws.Cells.Item(15, 6).Value = 58302
$ws.Cells.Item(16, 6).Value = 4976
$ws.Cells.Item(18, 6).Value = 538
$ws.Cells.Item(19, 6).Value = 122
$ws.Cells.Item(20, 6).Value = 946
$ws.Cells.Item(23, 6).Value = 5354
$ws.Cells.Item(25, 6).Value = 612
$ws.Cells.Item(26, 6).Value = 137
$ws.Cells.Item(28, 6).Value = 949
$ws.Cells.Item(29, 6).Value = 1462
$ws.Cells.Item(30, 6).Value = 2068
$ws.Cells.Item(32, 6).Value = 196
$ws.Cells.Item(33, 6).Value = 255
$ws.Cells.Item(36, 6).Value = 745
$ws.Cells.Item(39, 6).Value = 1198
$ws.Cells.Item(40, 6).Value = 418
$ws.Cells.Item(43, 6).Value = 227
$ws.Cells.Item(47, 6).Value = 2498

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 179
$ws.Cells.Item(4, 6).Value = 59
$ws.Cells.Item(5, 6).Value = 7747
$ws.Cells.Item(6, 6).Value = 129
$ws.Cells.Item(9, 6).Value = 4
$ws.Cells.Item(10, 6).Value = 15
$ws.Cells.Item(21, 6).Value = 45
$ws.Cells.Item(45, 6).Value = 47

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 2405
$ws.Cells.Item(5, 6).Value = 1632
$ws.Cells.Item(8, 6).Value = 2441
$ws.Cells.Item(9, 6).Value = 9477
$ws.Cells.Item(10, 6).Value = 1807
$ws.Cells.Item(11, 6).Value = 192
$ws.Cells.Item(12, 6).Value = 118
$ws.Cells.Item(15, 6).Value = 298
$ws.Cells.Item(16, 6).Value = 2560
$ws.Cells.Item(17, 6).Value = 269
$ws.Cells.Item(18, 6).Value = 84
$ws.Cells.Item(19, 6).Value = 558

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 358
$ws.Cells.Item(3, 6).Value = 3604
$ws.Cells.Item(4, 6).Value = 2405
$ws.Cells.Item(6, 6).Value = 1807
$ws.Cells.Item(7, 6).Value = 298
$ws.Cells.Item(8, 6).Value = 2560
$ws.Cells.Item(9, 6).Value = 269
$ws.Cells.Item(11, 6).Value = 7610
$ws.Cells.Item(12, 6).Value = 7824
$ws.Cells.Item(13, 6).Value = 4976
$ws.Cells.Item(14, 6).Value = 538
$ws.Cells.Item(15, 6).Value = 946
$ws.Cells.Item(17, 6).Value = 179
$ws.Cells.Item(18, 6).Value = 5354
$ws.Cells.Item(19, 6).Value = 612
$ws.Cells.Item(20, 6).Value = 137
$ws.Cells.Item(21, 6).Value = 84
$ws.Cells.Item(22, 6).Value = 1462
$ws.Cells.Item(23, 6).Value = 2068
$ws.Cells.Item(24, 6).Value = 129
$ws.Cells.Item(25, 6).Value = 558
$ws.Cells.Item(27, 6).Value = 15
$ws.Cells.Item(30, 6).Value = 196
$ws.Cells.Item(32, 6).Value = 745
$ws.Cells.Item(36, 6).Value = 45
$ws.Cells.Item(38, 6).Value = 418
$ws.Cells.Item(42, 6).Value = 227
$ws.Cells.Item(48, 6).Value = 2498
